# ---------------------------------------------------------------------------
# feat: add 2022-Q4 data
#
# 1. Insert a new worksheet "2022-Q4" right after "总计", pushing the other
#    quarterly sheets (2022-Q1, 2021-Q3, 2021-Q2, 2021-Q1, 2020-Q4) back by
#    one position.
# 2. Populate "2022-Q4" with its fund-holding detail rows.
# 3. Update the "总计" (summary) sheet: insert a new data row for 2022-Q4
#    right under the header, shifting the existing quarters down by a row,
#    and renumber the index column.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- locate sheets by their current (pre-edit) names -----------------------
$totalSheet = $wb.Worksheets.Item("总计")

# --- 1. insert the new "2022-Q4" sheet right after "总计" ------------------
$q4 = $wb.Worksheets.Add($null, $totalSheet)
$q4.Name = "2022-Q4"

# match page-setup of the rest of the workbook (boilerplate inches -> points)
$q4.PageSetup.LeftMargin = 54
$q4.PageSetup.RightMargin = 54
$q4.PageSetup.TopMargin = 72
$q4.PageSetup.BottomMargin = 72
$q4.PageSetup.HeaderMargin = 36
$q4.PageSetup.FooterMargin = 36

# --- 2. fill in the "2022-Q4" fund-holding detail table ---------------------
$q4.Cells.Item(1, 2).Value = "基金代码"
$q4.Cells.Item(1, 3).Value = "基金名称"
$q4.Cells.Item(1, 4).Value = "基金规模"
$q4.Cells.Item(1, 5).Value = "股票总仓位"
$q4.Cells.Item(1, 6).Value = "仓位占比"
$q4.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q4.Cells.Item(1, 8).Value = "仓位排名"

$q4.Cells.Item(2, 1).Value = 0

# the other quarterly sheets style their header row + index column with the
# same boxed/bold "s=2" format (pre-edit sheet2's header, now shifted to the
# "2022-Q1" sheet) -- clone it onto the new sheet's header/index cells.
$q1Sheet = $wb.Worksheets.Item("2022-Q1")
$q1Sheet.Cells.Item(1, 2).Copy() | Out-Null
$q4.Range("B1:H1").PasteSpecial(-4122)
$q1Sheet.Cells.Item(2, 1).Copy() | Out-Null
$q4.Cells.Item(2, 1).PasteSpecial(-4122)

$q4.Cells.Item(2, 2).Value = "'002597"
$q4.Cells.Item(2, 2).Style = "Normal"
$q4.Cells.Item(2, 3).Value = "'兴业成长动力灵活配置混合"
$q4.Cells.Item(2, 3).Style = "Normal"
$q4.Cells.Item(2, 4).Value = "'1.68"
$q4.Cells.Item(2, 4).Style = "Normal"
$q4.Cells.Item(2, 5).Value = "'89.03"
$q4.Cells.Item(2, 5).Style = "Normal"
$q4.Cells.Item(2, 6).Value = "'2.30"
$q4.Cells.Item(2, 6).Style = "Normal"
$q4.Cells.Item(2, 7).Value = "'0.0386"
$q4.Cells.Item(2, 7).Style = "Normal"
$q4.Cells.Item(2, 8).Value = 7

# --- 3. update the "总计" summary sheet -------------------------------------
# Existing rows 2..6 (2022-Q1 .. 2020-Q4) shift down to rows 3..7, then a new
# row 2 is written for 2022-Q4. Walk bottom-up so we never clobber a row
# before it has been read.
$dates  = @("2020-Q4", "2021-Q1", "2021-Q2", "2021-Q3", "2022-Q1")
$counts = @(5, 15, 7, 5, 2)
$values = @(1.37, 1.26, 0.52, 0.35, 0.04)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = 7 - $i
    $totalSheet.Cells.Item($row, 1).Value = 5 - $i
    $totalSheet.Cells.Item($row, 2).Value = $dates[$i]
    $totalSheet.Cells.Item($row, 3).Value = $counts[$i]
    $totalSheet.Cells.Item($row, 4).Value = $values[$i]
}

# new 2022-Q4 row (row 2)
$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 1
$totalSheet.Cells.Item(2, 4).Value = 0.04

# column-A cells use the bold/boxed "index" style throughout the data rows;
# copy that formatting from the header (B1) onto the newly written A7/A2.
$totalSheet.Cells.Item(1, 2).Copy() | Out-Null
$totalSheet.Cells.Item(7, 1).PasteSpecial(-4122)
$totalSheet.Cells.Item(1, 2).Copy() | Out-Null
$totalSheet.Cells.Item(2, 1).PasteSpecial(-4122)

# --- restore the originally-active tab ("2020-Q4", the last sheet) ---------
$wb.Worksheets.Item("2020-Q4").Activate()
